$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 49905
$ws.Range("J16").Value = 49905
$ws.Range("L16").Value = 49905
$ws.Range("N16").Value = -50365

$ws.Range("H55").Value = 354.96774
$ws.Range("I55").Value = 66.44444
$ws.Range("J55").Value = 473
$ws.Range("K55").Value = 66.44444
$ws.Range("L55").Value = 473
$ws.Range("M55").Value = 147.55556
$ws.Range("N55").Value = -901

$ws.Range("H127").Value = 1054.12
$ws.Range("I127").Value = 723.3939
$ws.Range("J127").Value = 1696.1177
$ws.Range("K127").Value = 2170.1817
$ws.Range("L127").Value = 5088.3531
$ws.Range("M127").Value = 2789.8183
$ws.Range("N127").Value = -15008.3531

$ws.Range("H131").Value = 957.5185
$ws.Range("I131").Value = 656.1667
$ws.Range("J131").Value = 3368.3333
$ws.Range("K131").Value = 1968.5001
$ws.Range("L131").Value = 10104.9999
$ws.Range("M131").Value = 3071.4999
$ws.Range("N131").Value = -20184.9999

$ws.Range("H132").Value = 1917340.5
$ws.Range("I132").Value = 2223357
$ws.Range("K132").Value = 6670071
$ws.Range("M132").Value = -6667541

$ws.Range("H138").Value = 3111.918
$ws.Range("I138").Value = 833.48486
$ws.Range("K138").Value = 2500.45458
$ws.Range("M138").Value = 2639.54542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1741.6875
$ws.Range("I2").Value = 1310.6666
$ws.Range("J2").Value = 2000.3
$ws.Range("K2").Value = 1310.6666
$ws.Range("L2").Value = 2000.3
$ws.Range("M2").Value = -1197.6666
$ws.Range("N2").Value = -2226.3

$ws.Range("H32").Value = 7529.729
$ws.Range("I32").Value = 3359.2896
$ws.Range("J32").Value = 23377.4
$ws.Range("K32").Value = 3359.2896
$ws.Range("L32").Value = 23377.4
$ws.Range("M32").Value = -3072.2896
$ws.Range("N32").Value = -23951.4

$ws.Range("H97").Value = 927.3684
$ws.Range("I97").Value = 756.9231
$ws.Range("J97").Value = 1296.6666
$ws.Range("K97").Value = 756.9231
$ws.Range("L97").Value = 1296.6666
$ws.Range("M97").Value = -260.9231
$ws.Range("N97").Value = -2288.6666

$ws.Range("H116").Value = 1741.6875
$ws.Range("I116").Value = 1310.6666
$ws.Range("J116").Value = 2000.3
$ws.Range("K116").Value = 1310.6666
$ws.Range("L116").Value = 2000.3
$ws.Range("M116").Value = 983.3334
$ws.Range("N116").Value = -6588.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1741.6875
$ws.Range("I3").Value = 1310.6666
$ws.Range("J3").Value = 2000.3
$ws.Range("K3").Value = 1310.6666
$ws.Range("L3").Value = 2000.3
$ws.Range("M3").Value = -1196.6666
$ws.Range("N3").Value = -2228.3

$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2226

$ws.Range("H20").Value = 41671216
$ws.Range("I20").Value = 83337980
$ws.Range("J20").Value = 4451.1665
$ws.Range("K20").Value = 83337980
$ws.Range("L20").Value = 4451.1665
$ws.Range("M20").Value = -83337733
$ws.Range("N20").Value = -4945.1665

$ws.Range("H94").Value = 977.7778
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -1702

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 22201
$ws.Range("I2").Value = 10250
$ws.Range("J2").Value = 70005
$ws.Range("K2").Value = 10250
$ws.Range("L2").Value = 70005
$ws.Range("M2").Value = -10137
$ws.Range("N2").Value = -70231

$ws.Range("H5").Value = 2581.4
$ws.Range("J5").Value = 2581.4
$ws.Range("L5").Value = 2581.4
$ws.Range("N5").Value = -2805.4

$ws.Range("H31").Value = 1940.5625
$ws.Range("I31").Value = 1658.5518
$ws.Range("J31").Value = 4666.6665
$ws.Range("K31").Value = 1658.5518
$ws.Range("L31").Value = 4666.6665
$ws.Range("M31").Value = -1363.5518
$ws.Range("N31").Value = -5256.6665

$ws.Range("H34").Value = 1940.5625
$ws.Range("I34").Value = 1658.5518
$ws.Range("J34").Value = 4666.6665
$ws.Range("K34").Value = 1658.5518
$ws.Range("L34").Value = 4666.6665
$ws.Range("M34").Value = -1456.5518
$ws.Range("N34").Value = -5070.6665

$ws.Range("H58").Value = 1749.3158
$ws.Range("I58").Value = 740
$ws.Range("J58").Value = 2109.7856
$ws.Range("K58").Value = 740
$ws.Range("L58").Value = 2109.7856
$ws.Range("M58").Value = -537
$ws.Range("N58").Value = -2515.7856

$ws.Range("H107").Value = 848.087
$ws.Range("I107").Value = 400.33334
$ws.Range("K107").Value = 400.33334
$ws.Range("M107").Value = 1519.66666

$ws.Range("H122").Value = 1224.381
$ws.Range("J122").Value = 1164.3636
$ws.Range("L122").Value = 3493.0908
$ws.Range("N122").Value = -8393.0908

$ws.Range("H136").Value = 1749.3158
$ws.Range("I136").Value = 740
$ws.Range("J136").Value = 2109.7856
$ws.Range("K136").Value = 2220
$ws.Range("L136").Value = 6329.3568
$ws.Range("M136").Value = 330
$ws.Range("N136").Value = -11429.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 252
$ws.Range("I98").Value = 128
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 384
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 1114
$ws.Range("N98").Value = -4496

$ws.Range("H113").Value = 427.19717
$ws.Range("I113").Value = 402.72726
$ws.Range("J113").Value = 448.44736
$ws.Range("K113").Value = 1208.18178
$ws.Range("L113").Value = 1345.34208
$ws.Range("M113").Value = 961.8182200000001
$ws.Range("N113").Value = -5685.34208

$ws.Range("H131").Value = 959.9434
$ws.Range("I131").Value = 481.66666
$ws.Range("J131").Value = 1021
$ws.Range("K131").Value = 1444.99998
$ws.Range("L131").Value = 3063
$ws.Range("M131").Value = 3595.00002
$ws.Range("N131").Value = -13143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3520.5454
$ws.Range("I102").Value = 2502.923
$ws.Range("J102").Value = 4990.4443
$ws.Range("K102").Value = 2502.923
$ws.Range("L102").Value = 4990.4443
$ws.Range("M102").Value = -880.9229999999998
$ws.Range("N102").Value = -8234.444299999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1328.7858
$ws.Range("I46").Value = 1375.125
$ws.Range("J46").Value = 1267
$ws.Range("K46").Value = 1375.125
$ws.Range("L46").Value = 1267
$ws.Range("M46").Value = -1187.125
$ws.Range("N46").Value = -1643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 32485.5
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 32485.5
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 32485.5
$ws.Range("N56").Value = -33913.5
$ws.Range("M56").ClearContents()

$ws.Range("H80").Value = 75041.57000000001
$ws.Range("J80").Value = 75041.57000000001
$ws.Range("L80").Value = 75041.57000000001
$ws.Range("N80").Value = -77037.57000000001

$ws.Range("H83").Value = 75041.57000000001
$ws.Range("J83").Value = 75041.57000000001
$ws.Range("L83").Value = 225124.71
$ws.Range("N83").Value = -235108.71

$ws.Range("H107").Value = 5492.0454
$ws.Range("I107").Value = 1016.8571
$ws.Range("J107").Value = 13323.625
$ws.Range("K107").Value = 3050.5713
$ws.Range("L107").Value = 39970.875
$ws.Range("M107").Value = -1130.5713
$ws.Range("N107").Value = -43810.875
